$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.322.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.869.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.26%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4675"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2843"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.36%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06536"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.82%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.38"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07873"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.870.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.090"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6747"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "277.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.309.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.497"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000007298"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.106.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.95%  "

$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.148"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "165.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.137"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.927"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.86%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.375"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09613"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.372"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.476"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.089"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04699"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.128"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7063"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("E37").Value = "  +0.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01853"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.283"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.532"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.945"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.90%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8493"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4174"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.46%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.75%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.149"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.190"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.61%  "

$ws.Range("E49").Value = "  -5.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.58%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05634"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.13%  "
